$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 491, shifting existing rows 491:586 down to 492:587
$ws.Rows(491).Insert()

# Populate the newly inserted row 491 with the new weekly price record
$ws.Range("A491").Value = 6
$ws.Range("B491").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C491").Value = "Metropolitana"
$ws.Range("D491").Value = 44508
$ws.Range("E491").Value = 13
$ws.Range("F491").Value = 100112028
$ws.Range("G491").Value = "Sandia"
$ws.Range("H491").Value = "Sin especificar"
$ws.Range("I491").Value = "Primera"
$ws.Range("J491").Value = 4500
$ws.Range("K491").Value = 650
$ws.Range("L491").Value = 700
$ws.Range("M491").Value = 671
$ws.Range("N491").Value = "$/kilo (volumen en unidades)"
$ws.Range("O491").Value = "Perú"
$ws.Range("P491").Value = 671
$ws.Range("Q491").Value = 1
$ws.Range("R491").Value = "Hortaliza"
